# Auto-generated edit script applying cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.293.50'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '2.900.80'
$ws.Range('E3').Value = '  +7.73%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '196.77'
$ws.Range('E5').Value = '  +4.65%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '601.31'
$ws.Range('E6').Value = '  +1.97%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.556'
$ws.Range('E8').Value = '  +2.74%  '
$ws.Range('E9').Value = '  -1.35%  '
$ws.Range('D10').Value = '2.899.74'
$ws.Range('E10').Value = '  +7.79%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.397'
$ws.Range('E11').Value = '  +10.88%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.160'
$ws.Range('E12').Value = '  -1.42%  '
$ws.Range('E13').Value = '  +4.71%  '
$ws.Range('D14').Value = '3.411.94'
$ws.Range('E14').Value = '  +6.86%  '
$ws.Range('D15').Value = '76.203.15'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.0000192'
$ws.Range('E16').Value = '  +1.48%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '27.64'
$ws.Range('E17').Value = '  +3.79%  '
$ws.Range('D18').Value = '2.900.79'
$ws.Range('E18').Value = '  +7.93%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '9.03'
$ws.Range('E19').Value = '  -3.46%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.63'
$ws.Range('E20').Value = '  +5.06%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '385.18'
$ws.Range('E21').Value = '  +2.74%  '
$ws.Range('E22').Value = '  +0.92%  '
$ws.Range('E23').Value = '  +2.39%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '72.00'
$ws.Range('E24').Value = '  +2.36%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.998'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('E26').Value = '  +2.36%  '
$ws.Range('E27').Value = '  +7.18%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.82'
$ws.Range('E28').Value = '  +4.21%  '
$ws.Range('E29').Value = '  +14.82%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.42'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '515.36'
$ws.Range('E32').Value = '  -1.45%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.85'
$ws.Range('E33').Value = '  +1.08%  '
$ws.Range('E34').Value = '  +3.54%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '165.32'
$ws.Range('E36').Value = '  +1.53%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '20.27'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.117'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '19.70'
$ws.Range('E39').Value = '  +1.55%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '184.25'
$ws.Range('E40').Value = '  +8.00%  '
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('E42').Value = '  +5.38%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '5.10'
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.69'
$ws.Range('E44').Value = '  -1.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0935'
$ws.Range('E45').Value = '  +10.54%  '
$ws.Range('E46').Value = '  +3.92%  '
$ws.Range('E47').Value = '  +1.90%  '
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.585'
$ws.Range('E49').Value = '  +8.57%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.690'
$ws.Range('E50').Value = '  +16.65%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '3.78'
$ws.Range('E51').Value = '  +3.43%  '

Write-Host "Applied 92 cell updates"
